$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = ''
$ws.Range("C8").Value = 94
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.0'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '0.00'

# Row 9
$ws.Range("C9").Value = 4
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Short point (up to 3 mtr.)'
$ws.Range("F9").Value = 256
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '1024.00'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'P. point'
$ws.Range("C10").Value = 15
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '4'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '9930.00'

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = 'P. point'
$ws.Range("C11").Value = 25
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'On board'
$ws.Range("F11").Value = 136
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '3400.00'

# Row 12
$ws.Range("C12").Value = 34
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.0'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 78
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '2652.00'

# Row 13
$ws.Range("C13").Value = 49
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.0'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 30
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1470.00'

# Row 14
$ws.Range("C14").Value = 86
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.0'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 219
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '18834.00'

# Row 15
$ws.Range("C15").Value = 66
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.0'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F15").Value = 303
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '19998.00'

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = ''
$ws.Range("C16").Value = 70
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.0'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F16").Value = 0
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '0.00'

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = 'R. mtr.'
$ws.Range("C17").Value = 37
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '17'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '25 mm'
$ws.Range("F17").Value = 56
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '2072.00'

# Row 18
$ws.Range("C18").Value = 66
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.0'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = 'Mtr.'
$ws.Range("C19").Value = 88
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F19").Value = 81
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '7128.00'

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = 'Mtr.'
$ws.Range("C20").Value = 7
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F20").Value = 122
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '854.00'

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = 'Set'
$ws.Range("C21").Value = 57
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.0'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F21").Value = 5733
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '326781.00'

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = ''
$ws.Range("C22").Value = 70
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.0'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F22").Value = 0
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '0.00'

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = 'Mtr.'
$ws.Range("C23").Value = 90
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F23").Value = 20
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '1800.00'

# Row 24
$ws.Range("C24").Value = 51
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.0'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'

# Row 25
$ws.Range("C25").Value = 39
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '16.0'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = ''
$ws.Range("C26").Value = 38
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.0'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F26").Value = 0
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '0.00'

# Row 27
$ws.Range("C27").Value = 62
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'

# Row 28
$ws.Range("C28").Value = 13
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '30'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = ' 6 A to 32 A rating'
$ws.Range("F28").Value = 187
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '2431.00'

# Row 29
$ws.Range("C29").Value = 47
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '31'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = 'Each'
$ws.Range("C30").Value = 60
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = ' 50/63 A rating'
$ws.Range("F30").Value = 900
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '54000.00'

# Row 31
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = ''
$ws.Range("C31").Value = 55
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '18.0'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F31").Value = 0
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '0.00'

# Row 32
$ws.Range("C32").Value = 56

# Row 34
$ws.Range("C34").Value = 10

# Row 35
$ws.Range("C35").Value = 78

# Row 36
$ws.Range("C36").Value = 81

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '500422.00'
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = '500422.00'

# Row 40
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '500422.00'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '500422.00'
